# testWorkbookClearCells.xlsx edit:
#  - drop the firstSheet view flag from the workbook window
#  - add a 3rd named region (region3) mirroring region1, placed at C18:I23
#    on the clearNamedRegion sheet
#  - populate that new region with the same header/data pattern as region1
#  - drop the extra "eight" column (I on region1 / M on region2) from both
#    existing regions on clearNamedRegion, which also prunes the now-unused
#    "eight" shared string
#  - update the saved selections on clearNamedRegion and clearRange

$wb = $excel.ActiveWorkbook

# --- workbook-level: new defined name -------------------------------------
$wb.Names.Add("region3", "=clearNamedRegion!`$C`$18:`$I`$23")

# --- clearNamedRegion: drop the 8th column of region1 (I2:I7) and
#     region2 (M10:M15) -----------------------------------------------------
$wsRegion = $wb.Worksheets.Item("clearNamedRegion")

$wsRegion.Range("I2:I7").ClearContents()
$wsRegion.Range("M10:M15").ClearContents()

# --- clearNamedRegion: write the new region3 block (C18:I23), a copy of
#     region1's header + data pattern shifted to start at C18 --------------
$wsRegion.Range("C18").Value = "one"
$wsRegion.Range("D18").Value = "two"
$wsRegion.Range("E18").Value = "three"
$wsRegion.Range("F18").Value = "four"
$wsRegion.Range("G18").Value = "five"
$wsRegion.Range("H18").Value = "six"
$wsRegion.Range("I18").Value = "seven"

for ($i = 0; $i -le 4; $i++) {
    $row = 19 + $i
    $wsRegion.Cells.Item($row, 3).Value = 1 + $i
    $wsRegion.Cells.Item($row, 4).Value = 6 + $i
    $wsRegion.Cells.Item($row, 5).Value = 11 + $i
    $wsRegion.Cells.Item($row, 6).Value = 16 + $i
    $wsRegion.Cells.Item($row, 7).Value = 21 + $i
    $wsRegion.Cells.Item($row, 8).Value = 26 + $i
    $wsRegion.Cells.Item($row, 9).Value = 31 + $i
}

# --- clearRange: move the saved selection to C3 ----------------------------
# (done before re-activating clearNamedRegion so that the final active tab
#  stays on clearNamedRegion, matching the workbook's activeTab)
$wsRange = $wb.Worksheets.Item("clearRange")
$wsRange.Range("C3").Select()

# --- clearNamedRegion: re-activate and move the saved selection to D11 ----
$wsRegion.Activate()
$wsRegion.Range("D11").Select()
